$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cryptocurrency data: refreshed prices, 1h volume
# percentages, and the reshuffled rows 48-51 (BabyDogeCoin dropped,
# RenderToken/Cronos/Algorand shifted up one row, EnergySwap added).
#
# Cells whose new text looks like a plain number (e.g. "215.23") need
# their number format forced to Text first, otherwise Excel silently
# reinterprets the assigned string as a floating point number.

$ws.Range('D2').Value = '27.031.18'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '1.653.10'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '215.23'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.65%  '
$ws.Range('E9').Value = '  +1.62%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '19.97'
$ws.Range('E10').Value = '  +4.28%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0868'
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('D12').Value = '1.886.32'
$ws.Range('E12').Value = '  +3.53%  '
$ws.Range('D13').Value = '1.665.71'
$ws.Range('E13').Value = '  +4.23%  '
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('E15').Value = '  +3.07%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '65.26'
$ws.Range('E16').Value = '  +2.82%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '239.54'
$ws.Range('E17').Value = '  +4.35%  '
$ws.Range('D18').Value = '27.032.41'
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('E19').Value = '  +2.04%  '
$ws.Range('E21').Value = '  +0.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.42'
$ws.Range('E22').Value = '  +4.23%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '2.24'
$ws.Range('E23').Value = '  +3.07%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '146.05'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('E26').Value = '  +0.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '7.12'
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('E29').Value = '  +2.93%  '
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('E32').Value = '  +3.19%  '
$ws.Range('D33').Value = '1.514.79'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('E35').Value = '  +8.75%  '
$ws.Range('E36').Value = '  -0.08%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.579'
$ws.Range('E37').Value = '  +1.46%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.890'
$ws.Range('E38').Value = '  +8.80%  '
$ws.Range('E39').Value = '  +3.20%  '
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  +4.15%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '65.80'
$ws.Range('E43').Value = '  +8.45%  '
$ws.Range('D44').Value = '1.793.66'
$ws.Range('E44').Value = '  +3.39%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.775'
$ws.Range('E45').Value = '  +2.26%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.916'
$ws.Range('E46').Value = '  -2.68%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '89.87'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.52'
$ws.Range('E48').Value = '  +2.81%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0508'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0976'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.56'
$ws.Range('E51').Value = '  +2.41%  '
